# "Add new RATE cards and update ID for Delta the Magnet Warrior"
#
# The RATE-JP sheet lists every card in the set. Column B holds each
# card's YGOPro id, which previously was only filled in for a handful of
# rows (one per named card). In this revision the id column was filled
# in for every row (a simple drag-fill continuing the existing
# 100911001, 100911002, ... sequence), three new cards were inserted
# into the name column (A7, A8 and A70), and the now-unused trailing
# blank template row (82) was removed since the sheet only needs 81
# rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RATE-JP")

# New card names that were added to the sheet (and to sharedStrings.xml).
$ws.Range("A7").Value  = "Speedroid Bamboo Horse"
$ws.Range("A8").Value  = "Wind Witch - Ice Bell"
$ws.Range("A70").Value = "Cipher Spectrum"

# Fill column B (card id) for every row from 3 to 81, continuing the
# sequence that already starts at row 2 with 100911001.
for ($row = 3; $row -le 81; $row++) {
    $ws.Cells.Item($row, 2).Value = 100911001 + ($row - 2)
}

# The sheet used to have a spare blank row (82) at the bottom; it's no
# longer needed now that every row has data, so remove it. This shifts
# the sheet dimension from A1:E82 down to A1:E81.
$ws.Rows.Item(82).Delete()

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("B9").Select()
